$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").Value = '''68.208.87'
    $ws.Range("E2").Value = '  -0.96%  '
    $ws.Range("D3").Value = '''3.790.30'
    $ws.Range("E3").Value = '  -2.02%  '
    $ws.Range("D4").Value = '''0.999'
    $ws.Range("E4").Value = '  -0.06%  '
    $ws.Range("D5").Value = '''597.55'
    $ws.Range("E5").Value = '  -0.85%  '
    $ws.Range("D6").Value = '''169.88'
    $ws.Range("E6").Value = '  -1.21%  '
    $ws.Range("D7").Value = '''3.789.90'
    $ws.Range("E7").Value = '  -2.01%  '
    $ws.Range("E8").Value = '  +0.02%  '
    $ws.Range("E9").Value = '  -0.48%  '
    $ws.Range("E10").Value = '  -3.22%  '
    $ws.Range("D11").Value = '''6.51'
    $ws.Range("E11").Value = '  +0.11%  '
    $ws.Range("D12").Value = '''0.454'
    $ws.Range("E12").Value = '  -1.97%  '
    $ws.Range("E13").Value = '  -2.96%  '
    $ws.Range("D14").Value = '''36.82'
    $ws.Range("E14").Value = '  -1.32%  '
    $ws.Range("D15").Value = '''4.430.45'
    $ws.Range("E15").Value = '  -2.04%  '
    $ws.Range("D16").Value = '''3.792.25'
    $ws.Range("E16").Value = '  -1.55%  '
    $ws.Range("D17").Value = '''18.73'
    $ws.Range("E17").Value = '  +1.16%  '
    $ws.Range("D18").Value = '''68.213.28'
    $ws.Range("E18").Value = '  -0.91%  '
    $ws.Range("D19").Value = '''7.20'
    $ws.Range("E19").Value = '  -2.99%  '
    $ws.Range("D20").Value = '''0.111'
    $ws.Range("E20").Value = '  -0.30%  '
    $ws.Range("E21").Value = '  -5.13%  '
    $ws.Range("D22").Value = '''468.16'
    $ws.Range("E22").Value = '  -1.15%  '
    $ws.Range("D23").Value = '''0.720'
    $ws.Range("E23").Value = '  -1.63%  '
    $ws.Range("E24").Value = '  -7.97%  '
    $ws.Range("D25").Value = '''83.86'
    $ws.Range("E25").Value = '  -0.06%  '
    $ws.Range("D26").Value = '''2.26'
    $ws.Range("E26").Value = '  -0.47%  '
    $ws.Range("D27").Value = '''12.18'
    $ws.Range("E27").Value = '  -0.09%  '
    $ws.Range("D28").Value = '''10.45'
    $ws.Range("E28").Value = '  -0.51%  '
    $ws.Range("E29").Value = '  -0.11%  '
    $ws.Range("D30").Value = '''3.942.50'
    $ws.Range("E30").Value = '  -2.06%  '
    $ws.Range("E31").Value = '  -1.30%  '
    $ws.Range("D32").Value = '''7.61'
    $ws.Range("E32").Value = '  -2.59%  '
    $ws.Range("D33").Value = '''30.62'
    $ws.Range("E33").Value = '  -2.58%  '
    $ws.Range("D34").Value = '''2.23'
    $ws.Range("E34").Value = '  -3.41%  '
    $ws.Range("D35").Value = '''9.24'
    $ws.Range("E35").Value = '  -1.74%  '
    $ws.Range("D36").Value = '''3.752.37'
    $ws.Range("E36").Value = '  -2.16%  '
    $ws.Range("D37").Value = '''3.76'
    $ws.Range("E37").Value = '  -5.98%  '
    $ws.Range("D38").Value = '''0.104'
    $ws.Range("E38").Value = '  -1.17%  '
    $ws.Range("D39").Value = '''0.139'
    $ws.Range("E39").Value = '  -0.94%  '
    $ws.Range("E40").Value = '  -1.86%  '
    $ws.Range("D41").Value = '''5.88'
    $ws.Range("E41").Value = '  -1.73%  '
    $ws.Range("D42").Value = '''1.00'
    $ws.Range("E42").Value = '  +0.01%  '
    $ws.Range("D43").Value = '''0.313'
    $ws.Range("E43").Value = '  -1.98%  '
    $ws.Range("D45").Value = '''8.69'
    $ws.Range("E45").Value = '  -0.52%  '
    $ws.Range("D46").Value = '''1.95'
    $ws.Range("E46").Value = '  -2.82%  '
    $ws.Range("D47").Value = '''406.04'
    $ws.Range("E47").Value = '  -4.28%  '
    $ws.Range("D48").Value = '''45.71'
    $ws.Range("E48").Value = '  -1.63%  '
    $ws.Range("D49").Value = '''0.000276'
    $ws.Range("E49").Value = '  -8.83%  '
    $ws.Range("E50").Value = '  +5.86%  '
    $ws.Range("D51").Value = '''143.06'
    $ws.Range("E51").Value = '  +0.32%  '
